$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 45866.33356521991

$ws.Range("A8").Value = 45866.37525955137
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

$ws.Range("B8").Value = 2025
$ws.Range("C8").Value = 31
$ws.Range("D8").Value = 15.72
$ws.Range("E8").Value = 86.02
$ws.Range("F8").Value = 147.66
$ws.Range("G8").Value = 8.41
$ws.Range("H8").Value = "ESE"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "09:00:22"
